$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Quoc V. Le), appearing first in Week 4, Lecture 3
$ws.Range("B27").Value = "Quoc V."
$ws.Range("A27").Value = "Le"
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3

# Update the selection to reflect where the user ended up after editing
$ws.Range("A29").Select()
